$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 1).Value = '3528 Ottoman Village Aged Care Broadmeadows'
$ws.Cells.Item(2, 2).Value = 24
$ws.Cells.Item(3, 1).Value = '3622 Olivet Care Aged Care Services Ringwood'
$ws.Cells.Item(3, 2).Value = 13
$ws.Cells.Item(4, 1).Value = '3652 Regis Aged Care Dandenong North'
$ws.Cells.Item(4, 2).Value = 25
$ws.Cells.Item(5, 1).Value = '3824 Estia Health South Morang'
$ws.Cells.Item(5, 2).Value = 57
$ws.Cells.Item(6, 1).Value = '3961 Heritage Water Gardens Aged Care Facility Sydenham'
$ws.Cells.Item(6, 2).Value = 17
$ws.Cells.Item(7, 1).Value = 'Aintree Primary School Aintree'
$ws.Cells.Item(7, 2).Value = 14
$ws.Cells.Item(8, 1).Value = 'Armstrong Creek School Armstrong Creek'
$ws.Cells.Item(8, 2).Value = 12
$ws.Cells.Item(9, 1).Value = 'Australian Meat Group Abattoir Dandenong South'
$ws.Cells.Item(9, 2).Value = 12
$ws.Cells.Item(10, 1).Value = 'Berwick Fields Primary School Berwick'
$ws.Cells.Item(10, 2).Value = 11
$ws.Cells.Item(11, 1).Value = 'Berwick Lodge Primary School Berwick'
$ws.Cells.Item(11, 2).Value = 22
$ws.Cells.Item(12, 1).Value = 'Bubup Womindjeka Family and Children''s Centre Port Melbourne'
$ws.Cells.Item(12, 2).Value = 11
$ws.Cells.Item(13, 1).Value = 'CREST Children''s Sanctuary Dandenong'
$ws.Cells.Item(13, 2).Value = 11
$ws.Cells.Item(14, 1).Value = 'Clifton Hill Primary School Clifton Hill'
$ws.Cells.Item(14, 2).Value = 13
$ws.Cells.Item(15, 1).Value = 'Dandenong North Primary School Dandenong'
$ws.Cells.Item(15, 2).Value = 10
$ws.Cells.Item(16, 1).Value = 'Elements Childcare Warralily Armstrong Creek'
$ws.Cells.Item(16, 2).Value = 21
$ws.Cells.Item(17, 1).Value = 'G & K OConnor PTY LTD Pakenham'
$ws.Cells.Item(17, 2).Value = 10
$ws.Cells.Item(18, 1).Value = 'Hamlyn Views School Hamlyn Heights'
$ws.Cells.Item(18, 2).Value = 10
$ws.Cells.Item(19, 1).Value = 'KingKids Early Learning Centre and Kindergarten Hallam'
$ws.Cells.Item(19, 2).Value = 11
$ws.Cells.Item(20, 1).Value = 'Lilydale Motor Inn Lilydale'
$ws.Cells.Item(20, 2).Value = 12
$ws.Cells.Item(21, 1).Value = 'Lowanna College Newborough'
$ws.Cells.Item(21, 2).Value = 35
$ws.Cells.Item(22, 1).Value = 'McQuinns Gym Bendigo'
$ws.Cells.Item(22, 2).Value = 14
$ws.Cells.Item(23, 1).Value = 'Metcash Limited Distribution Centre Laverton North'
$ws.Cells.Item(23, 2).Value = 12
$ws.Cells.Item(24, 1).Value = 'Monash Health Dandenong Hospital Dandenong'
$ws.Cells.Item(24, 2).Value = 10
$ws.Cells.Item(25, 1).Value = 'Morwell Park Primary School Morwell'
$ws.Cells.Item(25, 2).Value = 10
$ws.Cells.Item(26, 1).Value = 'Northern Bay College Wexford Campus Corio'
$ws.Cells.Item(26, 2).Value = 39
$ws.Cells.Item(27, 1).Value = 'Rosewood Downs Special Accommodation Home Dandenong'
$ws.Cells.Item(27, 2).Value = 13
$ws.Cells.Item(28, 1).Value = 'Saint Augustines Primary School Wodonga'
$ws.Cells.Item(28, 2).Value = 15
$ws.Cells.Item(29, 1).Value = 'Saint Monica''s Primary School Wodonga'
$ws.Cells.Item(29, 2).Value = 11
$ws.Cells.Item(30, 1).Value = 'St Brendans Primary School Shepparton'
$ws.Cells.Item(30, 2).Value = 12
$ws.Cells.Item(31, 1).Value = 'St Mary''s Primary School Swan Hill'
$ws.Cells.Item(31, 2).Value = 20
$ws.Cells.Item(32, 1).Value = 'St Thereses Primary School Kennington'
$ws.Cells.Item(32, 2).Value = 14
$ws.Cells.Item(33, 1).Value = 'St Vincents Hospital Emergency Department Melbourne'
$ws.Cells.Item(33, 2).Value = 18
$ws.Cells.Item(34, 1).Value = 'St. Brendans Catholic Primary School Lakes Entrance'
$ws.Cells.Item(34, 2).Value = 12
$ws.Cells.Item(35, 1).Value = 'TUROSI PTY LTD Thomastown'
$ws.Cells.Item(35, 2).Value = 14
$ws.Cells.Item(36, 1).Value = 'The Royal Children''s Hospital Parkville'
$ws.Cells.Item(36, 2).Value = 10
$ws.Cells.Item(37, 1).Value = 'Vizzarri Farms Koo Wee Rup'
$ws.Cells.Item(37, 2).Value = 14
$ws.Cells.Item(38, 1).Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Cells.Item(38, 2).Value = 42
$ws.Cells.Item(39, 1).Value = 'Western Health Sunshine Hospital Emergency Department St Albans'
$ws.Cells.Item(39, 2).Value = 14
$ws.Cells.Item(40, 1).Value = 'Wodonga Cemetery Wodonga'
$ws.Cells.Item(40, 2).Value = 39
$ws.Cells.Item(41, 1).Value = 'Wodonga Primary School Wodonga'
$ws.Cells.Item(41, 2).Value = 21
$ws.Cells.Item(42, 1).Value = 'Wodonga Senior Secondary College Wodonga'
$ws.Cells.Item(42, 2).Value = 20
$ws.Cells.Item(43, 1).Value = 'Wodonga South Primary School Wodonga'
$ws.Cells.Item(43, 2).Value = 34
$ws.Cells.Item(44, 1).Value = 'Woodend Primary School Woodend'
$ws.Cells.Item(44, 2).Value = 20
$ws.Cells.Item(45, 1).Value = 'Yallourn Power Station Yallourn'
$ws.Cells.Item(45, 2).Value = 10
$ws.Cells.Item(46, 1).Value = 'Yooralla Disability Residential Care Alfrieda Street St Albans'
$ws.Cells.Item(46, 2).Value = 12
